$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 1) was merged across the wrong columns (A1:B1) and its
# three header cells ("数据1","数据2","数据3") ended up shifted one column
# to the left relative to the detail rows below (row 2: 名字/描述/密码/邮箱,
# rows 3-4: the matching data). Fix: un-merge A1:B1, rotate every row's
# A:C values one column to the right (old C -> A, old A -> B, old B -> C),
# leaving column D untouched, then re-merge the group header across the
# now-correct columns (B1:C1).
#
# NOTE: values are moved with Range.Copy (a cut/paste-style move) rather
# than by reading .Value/.Value2 into a variable and writing it back --
# round-tripping a purely-numeric string like "123456" through .Value /
# .Value2 makes Excel re-interpret it as a number, which would change its
# stored type. Range.Copy preserves the original cell's stored type.
$rows = 1..4
$scratchCol = "Z"

# Remove the old (wrong) merge first so every cell in row 1 is individually
# addressable while we rotate the values around.
$ws.Range("A1:B1").UnMerge()

foreach ($r in $rows) {
    $a = "A$r"
    $b = "B$r"
    $c = "C$r"
    $z = "$scratchCol$r"

    # Stash old C so it isn't lost when B (or a clear) overwrites it.
    $ws.Range($c).Copy($ws.Range($z))

    if ($r -eq 1) {
        # Row 1's old B1 was blank (it was just the non-anchor half of the
        # A1:B1 merge), so the new C1 must end up blank too.
        $ws.Range($c).ClearContents()
    } else {
        $ws.Range($b).Copy($ws.Range($c))
    }

    $ws.Range($a).Copy($ws.Range($b))
    $ws.Range($z).Copy($ws.Range($a))
    $ws.Range($z).ClearContents()
}

# Re-create the merge over the correct columns for the "数据1" group header.
$ws.Range("B1:C1").Merge()
